# update weight excel file
#
# The "weight" sheet tracks a weekly weigh-in log (Thu..Wed across columns
# E:K). This commit fills in the rest of this week's entries (row 42) and
# starts next week's first entry (row 43), which also recalculates the
# "avg"/"std" helper columns (C/D) for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weight")
$ws.Activate()

# Row 42 (week 39) was missing Sun/Mon/Tue/Wed readings - fill them in.
$ws.Range("H42").Value = 61.1   # Sun
$ws.Range("I42").Value = 61.3   # Mon
$ws.Range("J42").Value = 61.1   # Tue
$ws.Range("K42").Value = 61     # Wed

# Row 43 (week 40) - record this week's first (Thursday) weigh-in.
$ws.Range("E43").Value = 61

# Reflect where the user ended up: scrolled a bit further down the frozen
# pane, with the active cell now on the new Thursday entry's row.
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$ws.Range("E44").Select()
